$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 904.6
$ws.Range("I28").Value = 904.6
$ws.Range("K28").Value = 904.6
$ws.Range("M28").Value = -419.6
$ws.Range("H31").Value = 2790.4285
$ws.Range("I31").Value = 2790.4285
$ws.Range("K31").Value = 8371.2855
$ws.Range("M31").Value = -8141.2855
$ws.Range("H32").Value = 18995.25
$ws.Range("I32").Value = 14588.8
$ws.Range("K32").Value = 14588.8
$ws.Range("M32").Value = -14262.8
$ws.Range("H41").Value = 1423.4286
$ws.Range("I41").Value = 1161.8572
$ws.Range("J41").Value = 1685
$ws.Range("K41").Value = 1161.8572
$ws.Range("L41").Value = 1685
$ws.Range("M41").Value = -721.8571999999999
$ws.Range("N41").Value = -2565
$ws.Range("H51").Value = 95999.60000000001
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 95999.60000000001
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 95999.60000000001
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -96967.60000000001
$ws.Range("H53").Value = 473.41666
$ws.Range("I53").Value = 209.6
$ws.Range("K53").Value = 209.6
$ws.Range("M53").Value = 427.4
$ws.Range("H81").Value = 47309.75
$ws.Range("J81").Value = 47309.75
$ws.Range("L81").Value = 47309.75
$ws.Range("N81").Value = -49305.75
$ws.Range("H84").Value = 47309.75
$ws.Range("J84").Value = 47309.75
$ws.Range("L84").Value = 141929.25
$ws.Range("N84").Value = -151913.25
$ws.Range("H106").Value = 2788
$ws.Range("I106").Value = 2796.6
$ws.Range("K106").Value = 2796.6
$ws.Range("M106").Value = -2165.6
$ws.Range("H107").Value = 1353.8889
$ws.Range("I107").Value = 1367.6364
$ws.Range("J107").Value = 1332.2858
$ws.Range("K107").Value = 1367.6364
$ws.Range("L107").Value = 1332.2858
$ws.Range("M107").Value = 552.3635999999999
$ws.Range("N107").Value = -5172.2858
$ws.Range("H113").Value = 100002130
$ws.Range("I113").Value = 33334878
$ws.Range("K113").Value = 33334878
$ws.Range("M113").Value = -33331624
$ws.Range("H135").Value = 7438.375
$ws.Range("J135").Value = 11493.846
$ws.Range("L135").Value = 103444.614
$ws.Range("N135").Value = -108514.614
$ws.Range("H141").Value = 2976.6667
$ws.Range("I141").Value = 2976.6667
$ws.Range("K141").Value = 8930.000100000001
$ws.Range("M141").Value = -3750.000100000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12196132
$ws.Range("I32").Value = 12500810
$ws.Range("K32").Value = 12500810
$ws.Range("M32").Value = -12500523
$ws.Range("H45").Value = 1960.6
$ws.Range("I45").Value = 887.8570999999999
$ws.Range("K45").Value = 887.8570999999999
$ws.Range("M45").Value = -510.8570999999999
$ws.Range("H52").Value = 115755.664
$ws.Range("J52").Value = 115755.664
$ws.Range("L52").Value = 115755.664
$ws.Range("N52").Value = -116391.664
$ws.Range("H61").Value = 41755840
$ws.Range("I61").Value = 100001200
$ws.Range("K61").Value = 100001200
$ws.Range("M61").Value = -100000988
$ws.Range("H136").Value = 41755840
$ws.Range("I136").Value = 100001200
$ws.Range("K136").Value = 300003600
$ws.Range("M136").Value = -300001050

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7366.091
$ws.Range("I20").Value = 6302
$ws.Range("J20").Value = 12154.5
$ws.Range("K20").Value = 6302
$ws.Range("L20").Value = 12154.5
$ws.Range("M20").Value = -6055
$ws.Range("N20").Value = -12648.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1062107.4
$ws.Range("I31").Value = 1546.9231
$ws.Range("J31").Value = 2594028
$ws.Range("K31").Value = 1546.9231
$ws.Range("L31").Value = 2594028
$ws.Range("M31").Value = -1251.9231
$ws.Range("N31").Value = -2594618
$ws.Range("H34").Value = 1062107.4
$ws.Range("I34").Value = 1546.9231
$ws.Range("J34").Value = 2594028
$ws.Range("K34").Value = 1546.9231
$ws.Range("L34").Value = 2594028
$ws.Range("M34").Value = -1344.9231
$ws.Range("N34").Value = -2594432
$ws.Range("H58").Value = 1781.5
$ws.Range("I58").Value = 1722.75
$ws.Range("J58").Value = 1899
$ws.Range("K58").Value = 1722.75
$ws.Range("L58").Value = 1899
$ws.Range("M58").Value = -1519.75
$ws.Range("N58").Value = -2305
$ws.Range("H122").Value = 5279.1577
$ws.Range("I122").Value = 2827.8333
$ws.Range("K122").Value = 8483.499899999999
$ws.Range("M122").Value = -6033.499899999999
$ws.Range("H134").Value = 772679.4
$ws.Range("I134").Value = 1001482.3
$ws.Range("K134").Value = 3004446.9
$ws.Range("M134").Value = -3001911.9
$ws.Range("H136").Value = 1781.5
$ws.Range("I136").Value = 1722.75
$ws.Range("J136").Value = 1899
$ws.Range("K136").Value = 5168.25
$ws.Range("L136").Value = 5697
$ws.Range("M136").Value = -2618.25
$ws.Range("N136").Value = -10797

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3959.8
$ws.Range("I70").Value = 899.5
$ws.Range("K70").Value = 2698.5
$ws.Range("M70").Value = -2383.5
$ws.Range("H73").Value = 3959.8
$ws.Range("I73").Value = 899.5
$ws.Range("K73").Value = 2698.5
$ws.Range("M73").Value = -1606.5
$ws.Range("H75").Value = 259912.5
$ws.Range("I75").Value = 2000
$ws.Range("J75").Value = 311495
$ws.Range("K75").Value = 6000
$ws.Range("L75").Value = 934485
$ws.Range("M75").Value = -5002
$ws.Range("N75").Value = -936481
$ws.Range("H78").Value = 259912.5
$ws.Range("I78").Value = 2000
$ws.Range("J78").Value = 311495
$ws.Range("K78").Value = 18000
$ws.Range("L78").Value = 2803455
$ws.Range("M78").Value = -13008
$ws.Range("N78").Value = -2813439
$ws.Range("H129").Value = 125089.375
$ws.Range("J129").Value = 199599
$ws.Range("L129").Value = 598797
$ws.Range("N129").Value = -608797

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 5000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H101").Value = 52972.832
$ws.Range("J101").Value = 52972.832
$ws.Range("L101").Value = 52972.832
$ws.Range("N101").Value = -59462.832
$ws.Range("H105").Value = 112827.5
$ws.Range("J105").Value = 112827.5
$ws.Range("L105").Value = 112827.5
$ws.Range("N105").Value = -119815.5
$ws.Range("H106").Value = 113619.664
$ws.Range("J106").Value = 113619.664
$ws.Range("L106").Value = 113619.664
$ws.Range("N106").Value = -116143.664
$ws.Range("H107").Value = 890.8570999999999
$ws.Range("I107").Value = 713.625
$ws.Range("J107").Value = 1127.1666
$ws.Range("K107").Value = 713.625
$ws.Range("L107").Value = 1127.1666
$ws.Range("M107").Value = 1206.375
$ws.Range("N107").Value = -4967.1666
$ws.Range("H123").Value = 39984
$ws.Range("J123").Value = 39984
$ws.Range("L123").Value = 39984
$ws.Range("N123").Value = -44884
$ws.Range("H126").Value = 6357.143
$ws.Range("I126").Value = 9750
$ws.Range("K126").Value = 29250
$ws.Range("M126").Value = -26780

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1249.25
$ws.Range("J16").Value = 556.3333
$ws.Range("L16").Value = 556.3333
$ws.Range("N16").Value = -896.3333
$ws.Range("H22").Value = 2067.0435
$ws.Range("I22").Value = 2168.3125
$ws.Range("J22").Value = 1835.5714
$ws.Range("K22").Value = 2168.3125
$ws.Range("L22").Value = 1835.5714
$ws.Range("M22").Value = -1873.3125
$ws.Range("N22").Value = -2425.5714
$ws.Range("H27").Value = 2067.0435
$ws.Range("I27").Value = 2168.3125
$ws.Range("J27").Value = 1835.5714
$ws.Range("K27").Value = 2168.3125
$ws.Range("L27").Value = 1835.5714
$ws.Range("M27").Value = -2061.3125
$ws.Range("N27").Value = -2049.5714
$ws.Range("H122").Value = 5137.727
$ws.Range("I122").Value = 4523.826
$ws.Range("K122").Value = 13571.478
$ws.Range("M122").Value = -11121.478
$ws.Range("H130").Value = 429
$ws.Range("J130").Value = 429
$ws.Range("L130").Value = 429
$ws.Range("N130").Value = -10469
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 38463256
$ws.Range("I107").Value = 55557660
$ws.Range("J107").Value = 848.25
$ws.Range("K107").Value = 166672980
$ws.Range("L107").Value = 2544.75
$ws.Range("M107").Value = -166671060
$ws.Range("N107").Value = -6384.75
